$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure all written cells keep their original Text cell type (source data
# is all inline/shared strings, even when the content looks numeric, e.g.
# "313.11" or "0.06390"). Setting NumberFormat to Text ("@") on each cell
# before assignment stops Excel from auto-converting the write into a
# Number cell and silently dropping things like trailing zeros.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.213.22"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.76%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.825.03"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.35%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.85%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.11"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.02%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.73%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4245"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.97%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3704"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.54%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07252"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8634"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.26%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.63%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.835.07"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.87%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.731"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.26%  "

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07094"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.28%  "

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.316"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.84%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.86"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.08%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.97%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008856"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.91%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.74%  "

$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.09"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.68%  "

$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.285.88"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.55%  "

$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.138"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.54%  "

$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.92"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.63%  "

$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.052.44"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.69%  "

$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.988"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.39%  "

$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.55"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.09%  "

$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.201"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.90%  "

$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.41"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.97%  "

$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.244"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.34%  "

$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.46"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.12%  "

$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08847"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.24%  "

$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "ARBITRUM"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.196"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.40%  "

$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7564"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.21%  "

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.457"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.78%  "

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "Frax"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.005"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.79%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.744"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -6.07%  "

$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.116"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.23%  "

$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01976"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.15%  "

$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05259"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.72%  "

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.337"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.19%  "

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.874"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.26%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1698"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.78%  "

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5056"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.54%  "

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.762"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.37%  "

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.63"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.93%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "107.38"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.06%  "

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4753"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.13%  "

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.004"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.75%  "

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06390"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.74%  "

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.673"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.26%  "

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.865"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.76%  "
